# Add a new "Dates" worksheet at the end of the workbook, populate it with
# date/time sample values formatted with a variety of date/time number
# formats, and make it the active sheet (mirrors the author's commit
# "first attempt to support parsing dates").

$wb = $excel.ActiveWorkbook

# --- Create the new worksheet after the last existing sheet ("Tab_entities") ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Dates"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 23
$ws.Columns.Item(4).ColumnWidth = 23.140625
$ws.Columns.Item(5).ColumnWidth = 28.28515625

# Helper-free, straight-line cell population (row by row, matching the sheet
# that Excel produced when the author started experimenting with date
# formats/styles).

# Row 1
$ws.Range("A1").Value = 44022
$ws.Range("A1").NumberFormat = "mm-dd-yy"

$ws.Range("B1").Value = 44022
$ws.Range("B1").NumberFormat = "dd/mm/yyyy;@"

$ws.Range("C1").NumberFormat = "@"
$ws.Range("C1").Value = "01.02.1789"
$ws.Range("C1").NumberFormat = "mm-dd-yy"

$ws.Range("D1").Value = 44022.122916666667
$ws.Range("D1").NumberFormat = "m/d/yy h:mm"

$ws.Range("E1").Value = 0.12359953703703704
$ws.Range("E1").NumberFormat = "h:mm:ss"

# Row 2
$ws.Range("A2").Value = 44022
$ws.Range("A2").NumberFormat = '[$-F800]dddd\,\ mmmm\ dd\,\ yyyy'

$ws.Range("B2").Value = 44022
$ws.Range("B2").NumberFormat = "dd/mm/yyyy;@"

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "31.12.1899"
$ws.Range("C2").NumberFormat = "mm-dd-yy"

$ws.Range("D2").Value = 44022.123599537037
$ws.Range("D2").NumberFormat = 'dd/mm/yyyy\ hh:mm:ss'

$ws.Range("E2").Value = 0.057638888888888885
$ws.Range("E2").NumberFormat = "[h]:mm:ss"

# Row 3
$ws.Range("A3").Value = 44022
$ws.Range("A3").NumberFormat = "dd/mm/yyyy;@"

$ws.Range("C3").Value = 1
$ws.Range("C3").NumberFormat = "mm-dd-yy"

$ws.Range("D3").NumberFormat = 'dd/mm/yyyy\ hh:mm:ss.000'

$ws.Range("E3").Value = 0.059938888888888889
$ws.Range("E3").NumberFormat = '[h]:mm:ss.000'

# Row 4
$ws.Range("A4").Value = 44022
$ws.Range("A4").NumberFormat = "dd/mm/yy;@"

$ws.Range("C4").Value = 2
$ws.Range("C4").NumberFormat = "mm-dd-yy"

# Row 5
$ws.Range("A5").Value = 44022
$ws.Range("A5").NumberFormat = 'dd/\ m/\ yy;@'

$ws.Range("C5").Value = 59
$ws.Range("C5").NumberFormat = "mm-dd-yy"

# Row 6
$ws.Range("A6").Value = 44022
$ws.Range("A6").NumberFormat = "d/m/yy;@"

$ws.Range("C6").Value = 60
$ws.Range("C6").NumberFormat = "mm-dd-yy"

# Row 7
$ws.Range("A7").Value = 44022
$ws.Range("A7").NumberFormat = 'yyyy\-mm\-dd;@'

$ws.Range("C7").Value = 61
$ws.Range("C7").NumberFormat = "mm-dd-yy"

# Row 8
$ws.Range("A8").Value = 44022
$ws.Range("A8").NumberFormat = '[$-100C]dddd\,\ d/\ mmmm\ yyyy;@'

$ws.Range("C8").Value = 929273
$ws.Range("C8").NumberFormat = "mm-dd-yy"

# Row 9
$ws.Range("A9").Value = 44022
$ws.Range("A9").NumberFormat = '[$-100C]d/\ mmmm\ yyyy;@'

# Row 10
$ws.Range("A10").Value = 44022
$ws.Range("A10").NumberFormat = '[$-100C]d\ mmm\ yy;@'

# Row 11
$ws.Range("A11").Value = 44022
$ws.Range("A11").NumberFormat = "mmm/yyyy"

# --- Page setup (matches paper size / orientation used by the author) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection / active cell mirrors the saved sheet state ---
$ws.Range("E15").Select()

Write-Output "Dates sheet added"
